$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Unicode subscript-three character used in one of the PEPE price strings
$sub3 = [char]0x2083

$ws.Range('D2').Value = '58.832.07'
$ws.Range('E2').Value = '  +1.06%  '
$ws.Range('D3').Value = '2.496.57'
$ws.Range('E3').Value = '  +1.89%  '
$ws.Range('E4').Value = '  +0.19%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '533.42'
$ws.Range('E5').Value = '  +1.21%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '136.19'
$ws.Range('E6').Value = '  +2.06%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.999'
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.564'
$ws.Range('E8').Value = '  +1.80%  '
$ws.Range('D9').Value = '2.509.62'
$ws.Range('E9').Value = '  +2.18%  '
$ws.Range('E10').Value = '  +2.77%  '
$ws.Range('E11').Value = '  -1.58%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '5.40'
$ws.Range('E12').Value = '  +2.21%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.348'
$ws.Range('E13').Value = '  +2.01%  '
$ws.Range('D14').Value = '2.944.87'
$ws.Range('E14').Value = '  +2.13%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '22.90'
$ws.Range('E15').Value = '  +1.85%  '
$ws.Range('D16').Value = '58.720.22'
$ws.Range('E16').Value = '  +1.10%  '
$ws.Range('E17').Value = '  +0.12%  '
$ws.Range('D18').Value = '2.505.97'
$ws.Range('E18').Value = '  +2.07%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '11.05'
$ws.Range('E19').Value = '  +3.54%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '4.25'
$ws.Range('E20').Value = '  +1.83%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '322.51'
$ws.Range('E21').Value = '  +0.81%  '
$ws.Range('B22').Value = 'Dai'
$ws.Range('C22').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.999'
$ws.Range('E22').Value = '  +0.30%  '
$ws.Range('B23').Value = 'Uniswap'
$ws.Range('C23').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '5.98'
$ws.Range('E23').Value = '  +4.79%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '65.24'
$ws.Range('E24').Value = '  +4.57%  '
$ws.Range('E25').Value = '  +3.56%  '
$ws.Range('E26').Value = '  -0.03%  '
$ws.Range('E27').Value = '  +1.60%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '7.50'
$ws.Range('E28').Value = '  +0.83%  '
$ws.Range('D29').Value = "0.0{0}0765" -f $sub3
$ws.Range('E29').Value = '  +2.62%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '6.55'
$ws.Range('E30').Value = '  +1.46%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '172.16'
$ws.Range('E31').Value = '  +5.41%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '1.74'
$ws.Range('E32').Value = '  +0.23%  '
$ws.Range('E33').Value = '  +8.61%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.999'
$ws.Range('E34').Value = '  +0.05%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '18.33'
$ws.Range('E35').Value = '  +0.91%  '
$ws.Range('E36').Value = '  +0.62%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '4.04'
$ws.Range('E37').Value = '  +0.94%  '
$ws.Range('E38').Value = '  -0.18%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '36.79'
$ws.Range('E39').Value = '  +1.31%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.806'
$ws.Range('E40').Value = '  +2.29%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '3.57'
$ws.Range('E41').Value = '  +1.44%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '283.15'
$ws.Range('E42').Value = '  +3.33%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '5.16'
$ws.Range('E43').Value = '  +3.02%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.996'
$ws.Range('E44').Value = '  -0.09%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.607'
$ws.Range('E45').Value = '  +3.77%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '130.21'
$ws.Range('E46').Value = '  +8.92%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '10.87'
$ws.Range('E47').Value = '  +0.28%  '
$ws.Range('E48').Value = '  +0.15%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.0502'
$ws.Range('E49').Value = '  -0.56%  '
$ws.Range('E50').Value = '  +0.76%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '17.26'
$ws.Range('E51').Value = '  +2.06%  '
